$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "reason" column (D), shifting it to E.
$ws.Columns.Item(4).Insert()

# New header for the inserted column; copy the bold/border/centered header
# formatting from the neighboring header cell, then set its text.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "project_count"

# Updated final_score values.
$ws.Range("C2").Value = 76.95
$ws.Range("C3").Value = 81.31999999999999
$ws.Range("C4").Value = 77.79000000000001

# New project_count values.
$ws.Range("D2").Value = 0.95
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1

# Updated reasoning text (now in column E).
$ws.Range("E2").Value = "`nReasoning: The candidate has demonstrated strong technical skills in the relevant technologies - ReactJS, JavaScript, CSS, Frontend Development, and NextJS - and has also successfully applied them in two projects. The candidate's skillset is a great fit for the job description and thus merits a score of 76.95."
$ws.Range("E3").Value = "`nReasoning: The candidate's projects demonstrate a solid understanding of the skills required for the job such as ReactJS, JavaScript, CSS, Frontend Development, NextJS, NodeJS, ExpressJS, Socket.IO, WebRTC, HTML, CSS, JS, Docker, Flutter, Dart, and Firebase. The projects also show a high level of proficiency in web design, user experience, web page optimization, and brand consistency, as well as communication and interpersonal skills. This is why the candidate has been given a score of 81.32."
$ws.Range("E4").Value = "`nReasoning: The candidate has demonstrated a strong understanding of the required skills (ReactJS, JavaScript, CSS, Frontend Development, NextJS) through their projects. The candidate has also used other technologies such as Pytorch, Tensorflow, Keras, Django Rest Framework, Python, and Sklearn to develop their projects, which shows their willingness to learn and adapt to new technologies. This, combined with their communication and interpersonal skills, makes them a suitable candidate for the Frontend Engineer Intern role."
